$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains weekly "Sandia" (watermelon) price records for the
# "Terminal La Palmera de La Serena" market. This commit adds the latest
# week's two new quality records (Extra / Super) at the top of the data
# block (row 41), pushing all the existing historical rows down by two
# rows (41-69 -> 43-71).

$ws.Rows("41:42").Insert()

# New row 41: Extra quality, week of 2021-12-24 (serial 44554)
$ws.Range("A41").Value = 8
$ws.Range("B41").Value = "Terminal La Palmera de La Serena"
$ws.Range("C41").Value = "Coquimbo"
$ws.Range("D41").Value = 44554
$ws.Range("E41").Value = 4
$ws.Range("F41").Value = 100112028
$ws.Range("G41").Value = "Sandia"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Extra"
$ws.Range("J41").Value = 1600
$ws.Range("K41").Value = 3500
$ws.Range("L41").Value = 4000
$ws.Range("M41").Value = 3750
$ws.Range("N41").Value = "`$/unidad"
$ws.Range("O41").Value = "Región de O'Higgins"
$ws.Range("P41").Value = 3750
$ws.Range("Q41").Value = 1
$ws.Range("R41").Value = "Hortaliza"

# New row 42: Super quality, week of 2021-12-24 (serial 44554)
$ws.Range("A42").Value = 8
$ws.Range("B42").Value = "Terminal La Palmera de La Serena"
$ws.Range("C42").Value = "Coquimbo"
$ws.Range("D42").Value = 44554
$ws.Range("E42").Value = 4
$ws.Range("F42").Value = 100112028
$ws.Range("G42").Value = "Sandia"
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("I42").Value = "Super"
$ws.Range("J42").Value = 2000
$ws.Range("K42").Value = 4500
$ws.Range("L42").Value = 5000
$ws.Range("M42").Value = 4750
$ws.Range("N42").Value = "`$/unidad"
$ws.Range("O42").Value = "Región de O'Higgins"
$ws.Range("P42").Value = 4750
$ws.Range("Q42").Value = 1
$ws.Range("R42").Value = "Hortaliza"
